# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 0
    6  = 1
    7  = 0
    8  = 1
    9  = 0
    10 = 1
    12 = 1
    13 = 1
    15 = 1
    16 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}

$wb.Save()
